$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "289.23"
Set-TextValue "E2" "0.79%"
Set-TextValue "E3" "1.88%"
Set-TextValue "D4" "5.263"
Set-TextValue "E4" "4.06%"
Set-TextValue "D5" "0.07077"
Set-TextValue "E5" "6.19%"
Set-TextValue "D6" "7.469"
Set-TextValue "E6" "1.70%"
Set-TextValue "E7" "5.14%"
Set-TextValue "D8" "1.394"
Set-TextValue "D9" "0.9053"
Set-TextValue "E9" "-3.88%"
Set-TextValue "E10" "3.52%"
Set-TextValue "D11" "0.07627"
Set-TextValue "E11" "15.70%"
Set-TextValue "D12" "0.07719"
Set-TextValue "E12" "1.59%"
Set-TextValue "D13" "0.02914"
Set-TextValue "E13" "-1.54%"
Set-TextValue "D14" "0.09022"
Set-TextValue "E14" "0.46%"
Set-TextValue "D15" "0.001594"
Set-TextValue "E15" "-1.24%"
Set-TextValue "D16" "0.0006535"
Set-TextValue "E16" "0.89%"
Set-TextValue "D17" "0.006149"
Set-TextValue "E17" "-3.34%"
Set-TextValue "E18" "1.19%"
Set-TextValue "D19" "2.233"
Set-TextValue "E19" "-0.85%"
Set-TextValue "D20" "0.3234"
Set-TextValue "E20" "0.60%"
Set-TextValue "D21" "0.1349"
Set-TextValue "E21" "2.85%"
Set-TextValue "D22" "4.001"
Set-TextValue "E22" "-2.09%"
Set-TextValue "D23" "0.1597"
Set-TextValue "E23" "2.61%"
Set-TextValue "D24" "0.04525"
Set-TextValue "E24" "0.64%"
Set-TextValue "D25" "0.001209"
Set-TextValue "E25" "2.09%"
Set-TextValue "D26" "0.004175"
Set-TextValue "E26" "-7.20%"
Set-TextValue "D27" "0.0001167"
Set-TextValue "E27" "-6.73%"
Set-TextValue "D28" "0.0001666"
Set-TextValue "E28" "2.82%"
Set-TextValue "D40" "0.04392"
Set-TextValue "E40" "4.49%"
Set-TextValue "D41" "0.007017"
Set-TextValue "E41" "3.96%"
Set-TextValue "D42" "0.1254"
Set-TextValue "E42" "-0.03%"
Set-TextValue "D43" "0.002065"
Set-TextValue "E43" "2.11%"
Set-TextValue "D44" "0.01194"
Set-TextValue "E44" "-2.82%"
Set-TextValue "D45" "0.00005828"
Set-TextValue "E45" "2.89%"
Set-TextValue "D47" "0.01297"
Set-TextValue "E47" "-0.86%"
